$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values (e.g. "1.00", "6.70")
# keep their exact original formatting instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '41.924.05'
$ws.Range("E2").Value = '  -4.60%  '
$ws.Range("D3").Value = '2.219.59'
$ws.Range("E3").Value = '  -5.65%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '244.89'
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  -6.32%  '
$ws.Range("D7").Value = '69.17'
$ws.Range("E7").Value = '  -6.54%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.547'
$ws.Range("E9").Value = '  -9.02%  '
$ws.Range("D10").Value = '0.0958'
$ws.Range("E10").Value = '  -5.37%  '
$ws.Range("D11").Value = '58.05'
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("D12").Value = '35.66'
$ws.Range("E12").Value = '  +7.24%  '
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("D14").Value = '6.70'
$ws.Range("E14").Value = '  -7.79%  '
$ws.Range("D15").Value = '2.549.31'
$ws.Range("E15").Value = '  -5.75%  '
$ws.Range("D16").Value = '14.80'
$ws.Range("E16").Value = '  -8.67%  '
$ws.Range("E17").Value = '  -7.42%  '
$ws.Range("D18").Value = '2.223.85'
$ws.Range("E18").Value = '  -5.62%  '
$ws.Range("D19").Value = '41.839.36'
$ws.Range("E19").Value = '  -4.70%  '
$ws.Range("D20").Value = '0.0₃0956'
$ws.Range("E20").Value = '  -7.13%  '
$ws.Range("D21").Value = '72.53'
$ws.Range("E21").Value = '  -7.17%  '
$ws.Range("D22").Value = '6.13'
$ws.Range("E22").Value = '  -7.82%  '
$ws.Range("D23").Value = '235.07'
$ws.Range("E23").Value = '  -7.10%  '
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  +9.65%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("E26").Value = '  -4.78%  '
$ws.Range("D27").Value = '2.46'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("D29").Value = '9.89'
$ws.Range("E29").Value = '  -5.28%  '
$ws.Range("D30").Value = '170.74'
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("D31").Value = '20.39'
$ws.Range("E31").Value = '  -8.38%  '
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  -5.52%  '
$ws.Range("E33").Value = '  -7.10%  '
$ws.Range("D34").Value = '0.0715'
$ws.Range("E34").Value = '  -4.40%  '
$ws.Range("D35").Value = '5.18'
$ws.Range("E35").Value = '  -3.82%  '
$ws.Range("D36").Value = '4.68'
$ws.Range("E36").Value = '  -7.79%  '
$ws.Range("D37").Value = '3.88'
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").Value = '22.78'
$ws.Range("E38").Value = '  +18.63%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '2.28'
$ws.Range("E39").Value = '  -4.88%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.0276'
$ws.Range("E40").Value = '  +1.06%  '
$ws.Range("D41").Value = '5.84'
$ws.Range("E41").Value = '  -9.11%  '
$ws.Range("D42").Value = '65.94'
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("D43").Value = '4.95'
$ws.Range("E43").Value = '  -9.69%  '
$ws.Range("D44").Value = '8.95'
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("E45").Value = '  -4.83%  '
$ws.Range("D46").Value = '0.190'
$ws.Range("E46").Value = '  -5.36%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = '4.52'
$ws.Range("E48").Value = '  +6.59%  '
$ws.Range("E49").Value = '  -3.70%  '
$ws.Range("D50").Value = '10.14'
$ws.Range("E50").Value = '  +7.93%  '
$ws.Range("D51").Value = '1.11'
$ws.Range("E51").Value = '  -4.39%  '
